# Add "SamplesTab" and "FilesTab" rows to the startup sheet (TabName/query/StatQuery/dbExcel/WebExcel),
# mirroring the existing CasesTab row. Commit: "Added Samples and Files Tab to all tests".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New Cypher query bodies (single-quoted here-strings => literal text, no `-escaping or `$`-expansion) ---
$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
 WHERE ss.disease_subtype IN ["Paget's Disease"] 
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
'@

$filesQuery = @'
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
 WHERE ss.disease_subtype IN ["Paget's Disease"] 
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
    order by f.file_name
'@

$neo4jFile = "TC10_Bento_Filter_Diagnosis-PagetsDisease_Neo4jData.xlsx"
$webFile   = "TC10_Bento_Filter_Diagnosis-PagetsDisease_WebData.xlsx"
$statQuery = $ws.Cells.Item(2, 3).Value()   # same StatQuery text reused by every tab row

# --- Column A tab names first, then column B query text (matches authoring order: both
#     tab labels typed in, then both long Cypher bodies pasted in) ---
$ws.Cells.Item(3, 1).Value = "SamplesTab"
$ws.Cells.Item(4, 1).Value = "FilesTab"

$ws.Cells.Item(3, 2).Value = $samplesQuery
$ws.Cells.Item(3, 2).WrapText = $true
$ws.Cells.Item(4, 2).Value = $filesQuery
$ws.Cells.Item(4, 2).WrapText = $true

# --- Row 3: SamplesTab ---
$ws.Cells.Item(3, 3).Value = $statQuery
$ws.Cells.Item(3, 3).WrapText = $true
$ws.Cells.Item(3, 4).Value = $neo4jFile
$ws.Cells.Item(3, 5).Value = $webFile
$ws.Rows.Item(3).RowHeight = 345.6

# --- Row 4: FilesTab ---
$ws.Cells.Item(4, 3).Value = $statQuery
$ws.Cells.Item(4, 3).WrapText = $true
$ws.Cells.Item(4, 4).Value = $neo4jFile
$ws.Cells.Item(4, 5).Value = $webFile
$ws.Rows.Item(4).RowHeight = 409.6

# --- Existing CasesTab row reflows slightly shorter once the sheet is resaved ---
$ws.Rows.Item(2).RowHeight = 316.8

# --- Columns re-fit to the new (wider/narrower) content ---
$ws.Columns.Item(1).ColumnWidth = 12.0
$ws.Columns.Item(2).ColumnWidth = 75.333333333333333
$ws.Columns.Item(3).ColumnWidth = 47.5
$ws.Columns.Item(4).ColumnWidth = 59.0
$ws.Columns.Item(5).ColumnWidth = 57.666666666666664

# --- Selection now spans the newly added SamplesTab/FilesTab rows ---
[void]$ws.Range("C2:E4").Select()

